# Weekly data update: insert a new week's record (2021-09-10, serial 44449)
# above the existing row 21, shifting rows 21-24 down to 22-25.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 21 (mirrors Excel's Rows(21).Insert behaviour:
# existing rows 21-24 shift down to 22-25, and formatting is inherited from
# the row above, which already carries the date style used in column D).
$ws.Rows.Item(21).Insert()

# Populate the newly inserted row 21 with the new weekly record.
$ws.Cells.Item(21, 1).Value = 9
$ws.Cells.Item(21, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(21, 3).Value = "Metropolitana"
$ws.Cells.Item(21, 4).Value = 44449
$ws.Cells.Item(21, 5).Value = 13
$ws.Cells.Item(21, 6).Value = 100112035
$ws.Cells.Item(21, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(21, 8).Value = "Sin especificar"
$ws.Cells.Item(21, 9).Value = "Primera"
$ws.Cells.Item(21, 10).Value = 18
$ws.Cells.Item(21, 11).Value = 24000
$ws.Cells.Item(21, 12).Value = 25000
$ws.Cells.Item(21, 13).Value = 24500
$ws.Cells.Item(21, 14).Value = "`$/malla 15 kilos"
$ws.Cells.Item(21, 15).Value = "Hijuelas"
$ws.Cells.Item(21, 16).Value = 1633
$ws.Cells.Item(21, 17).Value = 15
$ws.Cells.Item(21, 18).Value = "Hortaliza"
